$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.103.89"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.866.16"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.14"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4994"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3912"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09639"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +24.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.130"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +2.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.84"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.446"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.83"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.870.66"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.004"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.367"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001125"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.89"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06619"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.38"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.129"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.164.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.290"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.544"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.085.04"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.96%  "
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.02"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.46"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1059"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.053"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.601"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.38%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.614"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06741"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.432"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02385"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2176"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.991"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6260"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.004"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.68%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5974"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.670"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.272"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.34"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.982"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06828"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.31%  "
